# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Granada" (Vega Modelo de Temuco)
# at row 269, pushing the existing rows 269-305 down to 270-306.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 269 (shifts 269:305 -> 270:306).
$ws.Rows(269).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A269").Value = 10
$ws.Range("B269").Value = "Vega Modelo de Temuco"
$ws.Range("C269").Value = "La Araucanía"
$ws.Range("D269").Value = 45218
$ws.Range("E269").Value = 9
$ws.Range("F269").Value = "Fruta"
$ws.Range("G269").Value = 100104
$ws.Range("H269").Value = "Frutos de pepita"
$ws.Range("I269").Value = 100104001
$ws.Range("J269").Value = "Granada"
$ws.Range("K269").Value = "Wonderfull"
$ws.Range("L269").Value = "Primera"
$ws.Range("M269").Value = 180
$ws.Range("N269").Value = 17000
$ws.Range("O269").Value = 17000
$ws.Range("P269").Value = 17000
$ws.Range("Q269").Value = '$/bandeja 10 kilos granel'
$ws.Range("R269").Value = "Provincia de Limarí"
$ws.Range("S269").Value = 1700
$ws.Range("T269").Value = 10
